# Fruta / hortaliza, semanal
# Insert one new week of data (4 rows) above the most recent existing
# week in the "Naranja" sheet, pushing the rest of the table down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 4 blank rows at row 271 (old rows 271-313 become 275-317).
$ws.Range("A271:A274").EntireRow.Insert()

# Common / repeated values for this market block.
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limar" + [char]0x00ED
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100102
$producto  = "C" + [char]0x00ED + "tricos"
$categoriaId = 100102005
$categoria = "Naranja"
$unidad    = "`$/bins (400 kilos)"
$origen    = "Provincia de Limar" + [char]0x00ED
$kgUnidad  = 400

# New rows data:
# row, fecha, variedad, calidad, volumen, precioMin, precioMax, precioProm, precioKg
$newRows = @(
    @(271, 44476, "Lane Late",  "Primera", 20, 120000, 125000, 122500, 306),
    @(272, 44476, "Lane Late",  "Segunda", 20,  90000,  95000,  92500, 231),
    @(273, 44476, "Navel Late", "Primera", 20, 120000, 125000, 122500, 306),
    @(274, 44476, "Navel Late", "Segunda", 20,  90000,  95000,  92500, 231)
)

foreach ($r in $newRows) {
    $row        = $r[0]
    $fecha      = $r[1]
    $variedad   = $r[2]
    $calidad    = $r[3]
    $volumen    = $r[4]
    $precioMin  = $r[5]
    $precioMax  = $r[6]
    $precioProm = $r[7]
    $precioKg   = $r[8]

    $ws.Cells.Item($row, 1).Value2  = $mercadoId
    $ws.Cells.Item($row, 2).Value2  = $mercado
    $ws.Cells.Item($row, 3).Value2  = $region
    $ws.Cells.Item($row, 4).Value2  = $fecha
    $ws.Cells.Item($row, 5).Value2  = $codreg
    $ws.Cells.Item($row, 6).Value2  = $tipo
    $ws.Cells.Item($row, 7).Value2  = $productoId
    $ws.Cells.Item($row, 8).Value2  = $producto
    $ws.Cells.Item($row, 9).Value2  = $categoriaId
    $ws.Cells.Item($row, 10).Value2 = $categoria
    $ws.Cells.Item($row, 11).Value2 = $variedad
    $ws.Cells.Item($row, 12).Value2 = $calidad
    $ws.Cells.Item($row, 13).Value2 = $volumen
    $ws.Cells.Item($row, 14).Value2 = $precioMin
    $ws.Cells.Item($row, 15).Value2 = $precioMax
    $ws.Cells.Item($row, 16).Value2 = $precioProm
    $ws.Cells.Item($row, 17).Value2 = $unidad
    $ws.Cells.Item($row, 18).Value2 = $origen
    $ws.Cells.Item($row, 19).Value2 = $precioKg
    $ws.Cells.Item($row, 20).Value2 = $kgUnidad
}

Write-Output "Inserted new week rows 271-274; table now spans to row 317."
